# "Generate Report for Archive"
#  - Flip the localization status text from "Ready for handoff" to
#    "In Translation" everywhere it appears (Overview + per-language sheets).
#  - Narrow the language/status columns that previously showed that text.

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "In Translation"
$newWidth  = 13.4101845877511

foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    foreach ($cell in $used.Cells) {
        # Cast to [string] explicitly: some cells hold booleans (e.g. "True"),
        # and PowerShell's "-eq" coerces the string operand to the left-hand
        # side's type, which makes "$true -eq 'Ready for handoff'" true too.
        $val = [string]$cell.Value2
        if ($val -eq $oldStatus) {
            $cell.Value2 = $newStatus
        }
    }
}

# Overview sheet: the "zh-cn" (E) and "de-de" (F) status columns shrink.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Columns.Item(5).ColumnWidth = $newWidth
$wsOverview.Columns.Item(6).ColumnWidth = $newWidth

# zh-cn sheet: the "Status" (C) column shrinks.
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Columns.Item(3).ColumnWidth = $newWidth

# de-de sheet: the "Status" (C) column shrinks.
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Columns.Item(3).ColumnWidth = $newWidth
